$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.001.30"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "2.471.47"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.85"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.51"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "2.472.06"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.96"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").Value = "2.919.26"
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("D17").Value = "63.099.45"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "2.473.07"
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.19"
$ws.Range("E19").Value = "  +3.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.01"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "329.56"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.22"
$ws.Range("E22").Value = "  +8.10%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "663.35"
$ws.Range("E26").Value = "  +6.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.62"
$ws.Range("E27").Value = "  +14.22%  "
$ws.Range("D28").Value = "0.0₃0983"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "2.591.39"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("E30").Value = "  +811.61%  "
$ws.Range("E31").Value = "  +1.42%  "
$ws.Range("E32").Value = "  -1.64%  "
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.53"
$ws.Range("E35").Value = "  +2.77%  "
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.42"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "151.82"
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.71"
$ws.Range("E42").Value = "  -1.89%  "
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("E45").Value = "  +6.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "151.26"
$ws.Range("E46").Value = "  +4.47%  "
$ws.Range("E47").Value = "  +26.61%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.64"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0512"
$ws.Range("E51").Value = "  -1.10%  "
